# New weekly price record: insert a new row at position 150 (pushing the
# existing rows 150-193 down to 151-194) and populate it with this week's
# Durazno (peach) price entry for Vega Modelo de Temuco.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 150..193 down to 151..194, leaving a blank row 150 behind.
$ws.Rows.Item(150).Insert()

# Populate the newly inserted row 150 with the new weekly observation.
$ws.Cells.Item(150, 1).Value2  = 10
$ws.Cells.Item(150, 2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(150, 3).Value2  = "La Araucanía"
$ws.Cells.Item(150, 4).Value2  = 44559
$ws.Cells.Item(150, 5).Value2  = 9
$ws.Cells.Item(150, 6).Value2  = "Fruta"
$ws.Cells.Item(150, 7).Value2  = 100103
$ws.Cells.Item(150, 8).Value2  = "Frutos de hueso (carozo)"
$ws.Cells.Item(150, 9).Value2  = 100103004
$ws.Cells.Item(150, 10).Value2 = "Durazno"
$ws.Cells.Item(150, 11).Value2 = "Early Majestic"
$ws.Cells.Item(150, 12).Value2 = "Primera"
$ws.Cells.Item(150, 13).Value2 = 185
$ws.Cells.Item(150, 14).Value2 = 17000
$ws.Cells.Item(150, 15).Value2 = 18000
$ws.Cells.Item(150, 16).Value2 = 17595
$ws.Cells.Item(150, 17).Value2 = "$/bandeja 18 kilos granel"
$ws.Cells.Item(150, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(150, 19).Value2 = 978
$ws.Cells.Item(150, 20).Value2 = 18
